# Updated data to reflect new requirement separation.
# Split the "Prerequisites" column's trailing "Recommended: ..." clause
# into its own "Recommended" column, add empty "Corequisites" / "Concurrent"
# columns (all "NA"), and move "Terms Typically Offered" out to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the "Terms Typically Offered" column from D to G, and insert
#     the new Corequisites / Concurrent / Recommended columns in between.

# Header row
$ws.Range("G1").Value = $ws.Range("D1").Value()      # "Terms Typically Offered"
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Rows whose Prerequisites text has a trailing " Recommended: <text>" clause
# that needs to move into its own Recommended column.
$marker = " Recommended: "

for ($row = 2; $row -le 14; $row++) {
    $termsCell = $ws.Cells.Item($row, 4)   # column D, "Terms Typically Offered" (original)
    $terms = $termsCell.Value()

    $prereqCell = $ws.Cells.Item($row, 3)
    $prereq = $prereqCell.Value()
    $idx = $prereq.IndexOf($marker)

    if ($idx -ge 0) {
        $recommendedText = $prereq.Substring($idx + $marker.Length)
        $prereqCell.Value = $prereq.Substring(0, $idx)

        $ws.Cells.Item($row, 6).Value = $recommendedText     # column F (Recommended)
        $terms = $terms + " "
    } else {
        $ws.Cells.Item($row, 6).Value = "NA"                 # column F (Recommended)
    }

    $ws.Cells.Item($row, 4).Value = "NA"                     # column D (Corequisites)
    $ws.Cells.Item($row, 5).Value = "NA"                     # column E (Concurrent)
    $ws.Cells.Item($row, 7).Value = $terms                   # column G (Terms Typically Offered)
}
